$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1098.250371487543
$ws.Range("B3").Value = 1814.209095653003
$ws.Range("B4").Value = 2478.773449758434
$ws.Range("B5").Value = 2933.821426503184
$ws.Range("B6").Value = 3539.467590233888
$ws.Range("B7").Value = 4063.607084843637
$ws.Range("B8").Value = 4307.926644532105
$ws.Range("B9").Value = 5388.878248312813
$ws.Range("B10").Value = 5672.285609285031
$ws.Range("B11").Value = 6138.896130664856
$ws.Range("B12").Value = 6475.29518813123
$ws.Range("B13").Value = 6725.227449742923
$ws.Range("B14").Value = 7106.960042583044
$ws.Range("B15").Value = 7377.336667667126
$ws.Range("B16").Value = 7676.112476649897
$ws.Range("B17").Value = 7879.439636415213
$ws.Range("B18").Value = 8088.92646771216
$ws.Range("B19").Value = 8261.456869429348
$ws.Range("B20").Value = 7752.309965008676
$ws.Range("B21").Value = 7945.192260056127
$ws.Range("B22").Value = 8159.328026128601
$ws.Range("B23").Value = 8329.485314031115
$ws.Range("B24").Value = 8676.101320752452
$ws.Range("B25").Value = 8770.064407087555
$ws.Range("B26").Value = 8989.014387628475
$ws.Range("B27").Value = 9087.05865729228
$ws.Range("B28").Value = 9420.885093148152
$ws.Range("B29").Value = 9492.726952204463
$ws.Range("B30").Value = 9819.724336068961
$ws.Range("B31").Value = 9860.530968791916
$ws.Range("B32").Value = 9818.03404424881
$ws.Range("B33").Value = 9942.806041144024
$ws.Range("B34").Value = 10216.71801028242
$ws.Range("B35").Value = 10900.9068267242
$ws.Range("B36").Value = 10902.50100427957
$ws.Range("B37").Value = 10994.27361408337
$ws.Range("B38").Value = 11022.78091327783
$ws.Range("B39").Value = 11129.04893393894
$ws.Range("B40").Value = 11358.75241706292
$ws.Range("B41").Value = 11415.71444630126
$ws.Range("B42").Value = 11476.63772595368
$ws.Range("B43").Value = 11556.51183221017
$ws.Range("B44").Value = 11783.63582183473
$ws.Range("B45").Value = 11909.89417332255
$ws.Range("B46").Value = 11916.62386862896
$ws.Range("B47").Value = 12268.27636980571
$ws.Range("B48").Value = 12371.39972727417
$ws.Range("B49").Value = 12378.52041197571
$ws.Range("B50").Value = 12442.5913284068
$ws.Range("B51").Value = 12536.54493864434
$ws.Range("B52").Value = 12531.23186938742
$ws.Range("B53").Value = 12457.71218812057
$ws.Range("B54").Value = 12590.94598407466
$ws.Range("B55").Value = 12641.74294525916
$ws.Range("B56").Value = 12599.24780262792
$ws.Range("B57").Value = 12710.34760570354
$ws.Range("B58").Value = 12696.45595005546
$ws.Range("B59").Value = 12920.25523719204
$ws.Range("B60").Value = 12920.25523719204
$ws.Range("B61").Value = 12920.25523719204
$ws.Range("B62").Value = 12920.25523719204
